# Weekly data update: append Wk30 rows, refresh AutoFilter + view selection
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 147

# --- 1) (Re)apply the AutoFilter over the ORIGINAL data extent first -------
#        Doing this before the new rows are appended keeps the autoFilter's
#        <autoFilter ref=...> pinned to A1:H108 instead of auto-growing to
#        cover the freshly-appended rows below it.
$ws.Range("A1:H108").AutoFilter()
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "='Weekly Expenditure'!`$A`$1:`$H`$108")
$fdName.Visible = $false

# --- 2) Pre-format the new rows by cloning row 108's cell styles downward ---
$ws.Range("A108:H108").Copy()
$ws.Range("A109:H$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3) Write the Wk30 expenditure rows -------------------------------------
# Row 109
$ws.Range("A109").Value = 'Wk30'
$ws.Range("B109").Value = 'SAP'
$ws.Range("C109").Value = 'XS-PTS-0751'
$ws.Range("D109").Value = 'CDUJB8-15D Air Cylinder'
$ws.Range("E109").Value = '2pcs'
$ws.Range("F109").Value = 'Fishes'
$ws.Range("G109").Value = 45859
$ws.Range("H109").Value = 51.86
# Row 110
$ws.Range("A110").Value = 'Wk30'
$ws.Range("B110").Value = 'SAP'
$ws.Range("C110").Value = 11155497
$ws.Range("D110").Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Range("E110").Value = '20pcs'
$ws.Range("F110").Value = 'Fishes'
$ws.Range("G110").Value = 45859
$ws.Range("H110").Value = 702
# Row 111
$ws.Range("A111").Value = 'Wk30'
$ws.Range("B111").Value = 'SAP'
$ws.Range("C111").Value = 11155797
$ws.Range("D111").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E111").Value = '10pcs'
$ws.Range("F111").Value = 'Fishes'
$ws.Range("G111").Value = 45859
$ws.Range("H111").Value = 201.6
# Row 112
$ws.Range("A112").Value = 'Wk30'
$ws.Range("B112").Value = 'SAP'
$ws.Range("C112").Value = 11155497
$ws.Range("D112").Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Range("E112").Value = '20pcs'
$ws.Range("F112").Value = 'Fishes'
$ws.Range("G112").Value = 45860
$ws.Range("H112").Value = 702
# Row 113
$ws.Range("A113").Value = 'Wk30'
$ws.Range("B113").Value = 'SAP'
$ws.Range("C113").Value = 11155797
$ws.Range("D113").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E113").Value = '10pcs'
$ws.Range("F113").Value = 'Fishes'
$ws.Range("G113").Value = 45860
$ws.Range("H113").Value = 201.6
# Row 114
$ws.Range("A114").Value = 'Wk30'
$ws.Range("B114").Value = 'SAP'
$ws.Range("C114").Formula = '="11152231"'
$ws.Range("D114").Value = 'PTS-1090 TW.50.15.FI.0S.151.00 X2637'
$ws.Range("E114").Value = '4pcs'
$ws.Range("F114").Value = 'Fishes'
$ws.Range("G114").Value = 45860
$ws.Range("H114").Value = 2222.2399999999998
# Row 115
$ws.Range("A115").Value = 'Wk30'
$ws.Range("B115").Value = 'SAP'
$ws.Range("C115").Formula = '="11151246"'
$ws.Range("D115").Value = 'PTS-1069 TW.50.15.FI.0S.150.00 X2637'
$ws.Range("E115").Value = '4pcs'
$ws.Range("F115").Value = 'Fishes'
$ws.Range("G115").Value = 45860
$ws.Range("H115").Value = 2222.2399999999998
# Row 116
$ws.Range("A116").Value = 'Wk30'
$ws.Range("B116").Value = 'SAP'
$ws.Range("C116").Formula = '="11151236"'
$ws.Range("D116").Value = 'PTS-1059 70900.217 X2637 PRIMARY POGO'
$ws.Range("E116").Value = '400pcs'
$ws.Range("F116").Value = 'Fishes'
$ws.Range("G116").Value = 45860
$ws.Range("H116").Value = 840
# Row 117
$ws.Range("A117").Value = 'Wk30'
$ws.Range("B117").Value = 'SAP'
$ws.Range("C117").Value = 'XS-PTS-0356'
$ws.Range("D117").Value = 'GP-570D84A-03 Hyperspace Semicon Probes'
$ws.Range("E117").Value = '400pcs'
$ws.Range("F117").Value = 'Sihl'
$ws.Range("G117").Value = 45861
$ws.Range("H117").Value = 840
# Row 118
$ws.Range("A118").Value = 'Wk30'
$ws.Range("B118").Value = 'SAP'
$ws.Range("C118").Value = 'XS-PTS-0868'
$ws.Range("D118").Value = '6K-76235-H03X-ALN Hyperspace Semicon'
$ws.Range("E118").Value = '2pcs'
$ws.Range("F118").Value = 'Sihl'
$ws.Range("G118").Value = 45861
$ws.Range("H118").Value = 420.14
# Row 119
$ws.Range("A119").Value = 'Wk30'
$ws.Range("B119").Value = 'SAP'
$ws.Range("C119").Value = 'XS-PTS-0864'
$ws.Range("D119").Value = 'Model : 6K-76235-H03X-NST Hyperspace'
$ws.Range("E119").Value = '1pcs'
$ws.Range("F119").Value = 'Sihl'
$ws.Range("G119").Value = 45861
$ws.Range("H119").Value = 325.02999999999997
# Row 120
$ws.Range("A120").Value = 'Wk30'
$ws.Range("B120").Value = 'SAP'
$ws.Range("C120").Value = 'XS-PTS-0876'
$ws.Range("D120").Value = 'HX 2067 HPN CRC CO CONTACT CLEANER'
$ws.Range("E120").Value = '3pcs'
$ws.Range("F120").Value = 'Fishes'
$ws.Range("G120").Value = 45861
$ws.Range("H120").Value = 68.91
# Row 121
$ws.Range("A121").Value = 'Wk30'
$ws.Range("B121").Value = 'SAP'
$ws.Range("C121").Value = 11155497
$ws.Range("D121").Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Range("E121").Value = '16pcs'
$ws.Range("F121").Value = 'Fishes'
$ws.Range("G121").Value = 45861
$ws.Range("H121").Value = 561.6
# Row 122
$ws.Range("A122").Value = 'Wk30'
$ws.Range("B122").Value = 'SAP'
$ws.Range("C122").Value = 11155797
$ws.Range("D122").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E122").Value = '18pcs'
$ws.Range("F122").Value = 'Fishes'
$ws.Range("G122").Value = 45861
$ws.Range("H122").Value = 362.88
# Row 123
$ws.Range("A123").Value = 'Wk30'
$ws.Range("B123").Value = 'SAP'
$ws.Range("C123").Formula = '="11151237"'
$ws.Range("D123").Value = 'PTS-1060 70192.692 X2637 RUBBER TIP'
$ws.Range("E123").Value = '40pcs'
$ws.Range("F123").Value = 'Fishes'
$ws.Range("G123").Value = 45861
$ws.Range("H123").Value = 270.8
# Row 124
$ws.Range("A124").Value = 'Wk30'
$ws.Range("B124").Value = 'SAP'
$ws.Range("C124").Formula = '="11155143"'
$ws.Range("D124").Value = 'PTS-1136 300-001519-015 Semiconductor'
$ws.Range("E124").Value = '200pcs'
$ws.Range("F124").Value = 'Fishes'
$ws.Range("G124").Value = 45861
$ws.Range("H124").Value = 372
# Row 125
$ws.Range("A125").Value = 'Wk30'
$ws.Range("B125").Value = 'SAP'
$ws.Range("C125").Formula = '="11151250"'
$ws.Range("D125").Value = 'PTS-1073 10618.655 (10311.135){Q=129}'
$ws.Range("E125").Value = '1pcs'
$ws.Range("F125").Value = 'Fishes'
$ws.Range("G125").Value = 45861
$ws.Range("H125").Value = 1251.51
# Row 126
$ws.Range("A126").Value = 'Wk30'
$ws.Range("B126").Value = 'SAP'
$ws.Range("C126").Value = 'XS-SPM-0061'
$ws.Range("D126").Value = '14210873.140 SPAREP A218 - PICKUP'
$ws.Range("E126").Value = '1pcs'
$ws.Range("F126").Value = 'Lisa'
$ws.Range("G126").Value = 45861
$ws.Range("H126").Value = 207.12
# Row 127
$ws.Range("A127").Value = 'Wk30'
$ws.Range("B127").Value = 'SAP'
$ws.Range("C127").Value = 'XS-SPM-0006'
$ws.Range("D127").Value = '14210853.126 SPAREP A217 PICKUP MODULE'
$ws.Range("E127").Value = '1pcs'
$ws.Range("F127").Value = 'Lisa'
$ws.Range("G127").Value = 45861
$ws.Range("H127").Value = 220.49
# Row 128
$ws.Range("A128").Value = 'Wk30'
$ws.Range("B128").Value = 'SAP'
$ws.Range("C128").Value = 11155497
$ws.Range("D128").Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Range("E128").Value = '20pcs'
$ws.Range("F128").Value = 'Fishes'
$ws.Range("G128").Value = 45862
$ws.Range("H128").Value = 702
# Row 129
$ws.Range("A129").Value = 'Wk30'
$ws.Range("B129").Value = 'SAP'
$ws.Range("C129").Value = 11155797
$ws.Range("D129").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E129").Value = '10pcs'
$ws.Range("F129").Value = 'Fishes'
$ws.Range("G129").Value = 45862
$ws.Range("H129").Value = 201.6
# Row 130
$ws.Range("A130").Value = 'Wk30'
$ws.Range("B130").Value = 'SAP'
$ws.Range("C130").Formula = '="11151236"'
$ws.Range("D130").Value = 'PTS-1059 70900.217 X2637 PRIMARY POGO'
$ws.Range("E130").Value = '400pcs'
$ws.Range("F130").Value = 'Fishes'
$ws.Range("G130").Value = 45863
$ws.Range("H130").Value = 840
# Row 131
$ws.Range("A131").Value = 'Wk30'
$ws.Range("B131").Value = 'SAP'
$ws.Range("C131").Value = 11155497
$ws.Range("D131").Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Range("E131").Value = '16pcs'
$ws.Range("F131").Value = 'Fishes'
$ws.Range("G131").Value = 45863
$ws.Range("H131").Value = 561.6
# Row 132
$ws.Range("A132").Value = 'Wk30'
$ws.Range("B132").Value = 'SAP'
$ws.Range("C132").Value = 11155797
$ws.Range("D132").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E132").Value = '14pcs'
$ws.Range("F132").Value = 'Fishes'
$ws.Range("G132").Value = 45863
$ws.Range("H132").Value = 282.24
# Row 133
$ws.Range("A133").Value = 'Wk30'
$ws.Range("B133").Value = 'SAP'
$ws.Range("C133").Formula = '="11155143"'
$ws.Range("D133").Value = 'PTS-1136 300-001519-015 Semiconductor'
$ws.Range("E133").Value = '200pcs'
$ws.Range("F133").Value = 'Fishes'
$ws.Range("G133").Value = 45863
$ws.Range("H133").Value = 372
# Row 134
$ws.Range("A134").Value = 'Wk30'
$ws.Range("B134").Value = 'SAP'
$ws.Range("C134").Formula = '="11151245"'
$ws.Range("D134").Value = 'PTS-1068 TW.50.15.FI.0S.109.00 X2637_Stn'
$ws.Range("E134").Value = '1pcs'
$ws.Range("F134").Value = 'Fishes'
$ws.Range("G134").Value = 45863
$ws.Range("H134").Value = 304.08999999999997
# Row 135
$ws.Range("A135").Value = 'Wk30'
$ws.Range("B135").Value = 'SAP'
$ws.Range("C135").Value = 11151246
$ws.Range("D135").Value = 'PTS-1069 TW.50.15.FI.0S.150.00 X2637'
$ws.Range("E135").Value = '5pcs'
$ws.Range("F135").Value = 'Fishes'
$ws.Range("G135").Value = 45863
$ws.Range("H135").Value = 2777.8
# Row 136
$ws.Range("A136").Value = 'Wk30'
$ws.Range("B136").Value = 'SAP'
$ws.Range("C136").Value = 11155497
$ws.Range("D136").Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Range("E136").Value = '24pcs'
$ws.Range("F136").Value = 'Fishes'
$ws.Range("G136").Value = 45863
$ws.Range("H136").Value = 842.4
# Row 137
$ws.Range("A137").Value = 'Wk30'
$ws.Range("B137").Value = 'SAP'
$ws.Range("C137").Value = 11155797
$ws.Range("D137").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E137").Value = '16pcs'
$ws.Range("F137").Value = 'Fishes'
$ws.Range("G137").Value = 45863
$ws.Range("H137").Value = 322.56
# Row 138
$ws.Range("A138").Value = 'Wk30'
$ws.Range("B138").Value = 'SAP'
$ws.Range("C138").Formula = '="11154654"'
$ws.Range("D138").Value = 'PTS-1097 10618.836 (10412.950) X2544'
$ws.Range("E138").Value = '1pcs'
$ws.Range("F138").Value = 'Fishes'
$ws.Range("G138").Value = 45863
$ws.Range("H138").Value = 3018.87
# Row 139
$ws.Range("A139").Value = 'Wk30'
$ws.Range("B139").Value = 'SAP'
$ws.Range("C139").Value = 11155797
$ws.Range("D139").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E139").Value = '10pcs'
$ws.Range("F139").Value = 'Fishes'
$ws.Range("G139").Value = 45864
$ws.Range("H139").Value = 201.6
# Row 140
$ws.Range("A140").Value = 'Wk30'
$ws.Range("B140").Value = 'SAP'
$ws.Range("C140").Formula = '="11151255"'
$ws.Range("D140").Value = 'PTS-1078 TW.50.1A.FI.0S.125.01 TOP PIN'
$ws.Range("E140").Value = '1pcs'
$ws.Range("F140").Value = 'Fishes'
$ws.Range("G140").Value = 45864
$ws.Range("H140").Value = 877.19
# Row 141
$ws.Range("A141").Value = 'Wk30'
$ws.Range("B141").Value = 'SAP'
$ws.Range("C141").Value = 11155497
$ws.Range("D141").Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Range("E141").Value = '10pcs'
$ws.Range("F141").Value = 'Fishes'
$ws.Range("G141").Value = 45864
$ws.Range("H141").Value = 351
# Row 142
$ws.Range("A142").Value = 'Wk30'
$ws.Range("B142").Value = 'SAP'
$ws.Range("C142").Value = 11155797
$ws.Range("D142").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E142").Value = '8pcs'
$ws.Range("F142").Value = 'Fishes'
$ws.Range("G142").Value = 45864
$ws.Range("H142").Value = 161.28
# Row 143
$ws.Range("A143").Value = 'Wk30'
$ws.Range("B143").Value = 'SAP'
$ws.Range("C143").Value = 11152231
$ws.Range("D143").Value = 'PTS-1090 TW.50.15.FI.0S.151.00 X2637'
$ws.Range("E143").Value = '1pcs'
$ws.Range("F143").Value = 'Fishes'
$ws.Range("G143").Value = 45864
$ws.Range("H143").Value = 555.55999999999995
# Row 144
$ws.Range("A144").Value = 'Wk30'
$ws.Range("B144").Value = 'SAP'
$ws.Range("C144").Value = 11151246
$ws.Range("D144").Value = 'PTS-1069 TW.50.15.FI.0S.150.00 X2637'
$ws.Range("E144").Value = '1pcs'
$ws.Range("F144").Value = 'Fishes'
$ws.Range("G144").Value = 45864
$ws.Range("H144").Value = 555.55999999999995
# Row 145
$ws.Range("A145").Value = 'Wk30'
$ws.Range("B145").Value = 'SAP'
$ws.Range("C145").Formula = '="11154654"'
$ws.Range("D145").Value = 'PTS-1097 10618.836 (10412.950){Q=16},'
$ws.Range("E145").Value = '1pcs'
$ws.Range("F145").Value = 'Fishes'
$ws.Range("G145").Value = 45864
$ws.Range("H145").Value = 3018.87
# Row 146
$ws.Range("A146").Value = 'Wk30'
$ws.Range("B146").Value = 'SAP'
$ws.Range("C146").Value = 11155497
$ws.Range("D146").Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Range("E146").Value = '16pcs'
$ws.Range("F146").Value = 'Fishes'
$ws.Range("G146").Value = 45865
$ws.Range("H146").Value = 561.6
# Row 147
$ws.Range("A147").Value = 'Wk30'
$ws.Range("B147").Value = 'SAP'
$ws.Range("C147").Value = 11155797
$ws.Range("D147").Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Range("E147").Value = '10pcs'
$ws.Range("F147").Value = 'Fishes'
$ws.Range("G147").Value = 45865
$ws.Range("H147").Value = 201.6

# --- 4) Collapse the helper formulas used to force text-typed numeric codes -
#        (keeps them as plain shared-string cells, matching manually-typed data)
$forceTextCells = @("C114","C115","C116","C123","C124","C125","C130","C133","C134","C138","C140","C145")
foreach ($addr in $forceTextCells) {
    $cell = $ws.Range($addr)
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

# --- 5) Reset the view: scroll to top, select A2 ----------------------------
$ws.Range("A2").Select()

Write-Host "Wk30 rows appended through row $lastRow"
